# Rewrite of example 11 unit tests
# Updates the "ATDD Scenarios" sheet in the LookupValue - Chapter 12 workbook:
# replaces event-subscriber-oriented scenario/when text with the new
# CheckLookupvalueExistsOnSalesHeader / InheritLookupValue* / ApplyLookupValue*
# naming, and adjusts the affected (wrap-text) row heights to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATDD Scenarios")

# --- Scenario 0100: Check failure CheckLookupvalueExistsOnSalesHeader Sales Posting ---
$ws.Range("F3").Value = "Check failure CheckLookupvalueExistsOnSalesHeader Sales Posting"
$ws.Range("H5").Value = "Trigger CheckLookupvalueExistsOnSalesHeader Sales Posting"

# --- Scenario 0101: Check success CheckLookupvalueExistsOnSalesHeader Sales Posting ---
$ws.Range("F7").Value = "Check success CheckLookupvalueExistsOnSalesHeader Sales Posting"
$ws.Range("H9").Value = "Trigger CheckLookupvalueExistsOnSalesHeader Sales Posting"

# --- Scenario 0102: Check failure CheckLookupvalueExistsOnSalesHeader Whse. Posting ---
$ws.Range("F11").Value = "Check failure CheckLookupvalueExistsOnSalesHeader Whse. Posting"
$ws.Range("H13").Value = "Trigger CheckLookupvalueExistsOnSalesHeader Whse. Posting"

# --- Scenario 0103: Check success CheckLookupvalueExistsOnSalesHeader Whse. Posting ---
$ws.Range("F15").Value = "Check success CheckLookupvalueExistsOnSalesHeader Whse. Posting"
$ws.Range("H17").Value = "Trigger CheckLookupvalueExistsOnSalesHeader Whse. Posting"

# --- Scenario 0104: Check InheritLookupValueFromCustomer ---
$ws.Range("F20").Value = "Check InheritLookupValueFromCustomer"
$ws.Range("H23").Value = "Trigger InheritLookupValueFromCustomer"

# --- Scenario 0105: Check ApplyLookupValueFromCustomerTemplate from Contact ---
$ws.Range("F25").Value = "Check ApplyLookupValueFromCustomerTemplate from Contact"
$ws.Range("H28").Value = "Trigger ApplyLookupValueFromCustomerTemplate"

# --- Scenario 0106: Check ApplyLookupValueFromCustomerTemplate ---
$ws.Range("F30").Value = "Check ApplyLookupValueFromCustomerTemplate"
$ws.Range("H33").Value = "Trigger ApplyLookupValueFromCustomerTemplate"

# --- Scenario 0107: Check InheritLookupValueFromSalesHeader ---
$ws.Range("F36").Value = "Check InheritLookupValueFromSalesHeader"
$ws.Range("H39").Value = "Trigger InheritLookupValueFromSalesHeader"

# --- Row heights (wrap-text auto-fit in response to the new text lengths) ---
$ws.Rows.Item(3).RowHeight = 45.75
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30.75
$ws.Rows.Item(23).RowHeight = 16.5
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(30).RowHeight = 30
$ws.Rows.Item(36).RowHeight = 30.75

# --- Reset the active selection back to A1 (was A15) ---
$ws.Range("A1").Select()
